# Append rows 175-209 to the "Teste" sheet: a tweet (shared text that already
# exists elsewhere in this workbook) plus its relevance label (0/1).
$wb = $excel.ActiveWorkbook
$wsTreinamento = $wb.Worksheets.Item("Treinamento")
$wsTeste = $wb.Worksheets.Item("Teste")

$text = $wsTreinamento.Cells.Item(232, 1).Value2
$wsTeste.Cells.Item(175, 1).Value = $text
$wsTeste.Cells.Item(175, 2).Value = 0
$wsTeste.Rows.Item(175).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(128, 1).Value2
$wsTeste.Cells.Item(176, 1).Value = $text
$wsTeste.Cells.Item(176, 2).Value = 0
$wsTeste.Rows.Item(176).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(61, 1).Value2
$wsTeste.Cells.Item(177, 1).Value = $text
$wsTeste.Cells.Item(177, 2).Value = 1
$wsTeste.Rows.Item(177).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(278, 1).Value2
$wsTeste.Cells.Item(178, 1).Value = $text
$wsTeste.Cells.Item(178, 2).Value = 0
$wsTeste.Rows.Item(178).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(179, 1).Value2
$wsTeste.Cells.Item(179, 1).Value = $text
$wsTeste.Cells.Item(179, 2).Value = 0
$wsTeste.Rows.Item(179).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(122, 1).Value2
$wsTeste.Cells.Item(180, 1).Value = $text
$wsTeste.Cells.Item(180, 2).Value = 1
$wsTeste.Rows.Item(180).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(105, 1).Value2
$wsTeste.Cells.Item(181, 1).Value = $text
$wsTeste.Cells.Item(181, 2).Value = 0
$wsTeste.Rows.Item(181).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(64, 1).Value2
$wsTeste.Cells.Item(182, 1).Value = $text
$wsTeste.Cells.Item(182, 2).Value = 1
$wsTeste.Rows.Item(182).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(111, 1).Value2
$wsTeste.Cells.Item(183, 1).Value = $text
$wsTeste.Cells.Item(183, 2).Value = 0
$wsTeste.Rows.Item(183).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(121, 1).Value2
$wsTeste.Cells.Item(184, 1).Value = $text
$wsTeste.Cells.Item(184, 2).Value = 0
$wsTeste.Rows.Item(184).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(279, 1).Value2
$wsTeste.Cells.Item(185, 1).Value = $text
$wsTeste.Cells.Item(185, 2).Value = 1
$wsTeste.Rows.Item(185).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(31, 1).Value2
$wsTeste.Cells.Item(186, 1).Value = $text
$wsTeste.Cells.Item(186, 2).Value = 0
$wsTeste.Rows.Item(186).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(156, 1).Value2
$wsTeste.Cells.Item(187, 1).Value = $text
$wsTeste.Cells.Item(187, 2).Value = 1
$wsTeste.Rows.Item(187).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(67, 1).Value2
$wsTeste.Cells.Item(188, 1).Value = $text
$wsTeste.Cells.Item(188, 2).Value = 1
$wsTeste.Rows.Item(188).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(27, 1).Value2
$wsTeste.Cells.Item(189, 1).Value = $text
$wsTeste.Cells.Item(189, 2).Value = 0
$wsTeste.Rows.Item(189).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(25, 1).Value2
$wsTeste.Cells.Item(190, 1).Value = $text
$wsTeste.Cells.Item(190, 2).Value = 1
$wsTeste.Rows.Item(190).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(72, 1).Value2
$wsTeste.Cells.Item(191, 1).Value = $text
$wsTeste.Cells.Item(191, 2).Value = 1
$wsTeste.Rows.Item(191).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(148, 1).Value2
$wsTeste.Cells.Item(192, 1).Value = $text
$wsTeste.Cells.Item(192, 2).Value = 1
$wsTeste.Rows.Item(192).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(30, 1).Value2
$wsTeste.Cells.Item(193, 1).Value = $text
$wsTeste.Cells.Item(193, 2).Value = 1
$wsTeste.Rows.Item(193).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(103, 1).Value2
$wsTeste.Cells.Item(194, 1).Value = $text
$wsTeste.Cells.Item(194, 2).Value = 0
$wsTeste.Rows.Item(194).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(201, 1).Value2
$wsTeste.Cells.Item(195, 1).Value = $text
$wsTeste.Cells.Item(195, 2).Value = 1
$wsTeste.Rows.Item(195).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(100, 1).Value2
$wsTeste.Cells.Item(196, 1).Value = $text
$wsTeste.Cells.Item(196, 2).Value = 1
$wsTeste.Rows.Item(196).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(35, 1).Value2
$wsTeste.Cells.Item(197, 1).Value = $text
$wsTeste.Cells.Item(197, 2).Value = 1
$wsTeste.Rows.Item(197).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(295, 1).Value2
$wsTeste.Cells.Item(198, 1).Value = $text
$wsTeste.Cells.Item(198, 2).Value = 1
$wsTeste.Rows.Item(198).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(91, 1).Value2
$wsTeste.Cells.Item(199, 1).Value = $text
$wsTeste.Cells.Item(199, 2).Value = 0
$wsTeste.Rows.Item(199).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(77, 1).Value2
$wsTeste.Cells.Item(200, 1).Value = $text
$wsTeste.Cells.Item(200, 2).Value = 1
$wsTeste.Rows.Item(200).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(108, 1).Value2
$wsTeste.Cells.Item(201, 1).Value = $text
$wsTeste.Cells.Item(201, 2).Value = 0
$wsTeste.Rows.Item(201).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(153, 1).Value2
$wsTeste.Cells.Item(202, 1).Value = $text
$wsTeste.Cells.Item(202, 2).Value = 1
$wsTeste.Rows.Item(202).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(95, 1).Value2
$wsTeste.Cells.Item(203, 1).Value = $text
$wsTeste.Cells.Item(203, 2).Value = 1
$wsTeste.Rows.Item(203).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(8, 1).Value2
$wsTeste.Cells.Item(204, 1).Value = $text
$wsTeste.Cells.Item(204, 2).Value = 1
$wsTeste.Rows.Item(204).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(268, 1).Value2
$wsTeste.Cells.Item(205, 1).Value = $text
$wsTeste.Cells.Item(205, 2).Value = 1
$wsTeste.Rows.Item(205).EntireRow.AutoFit()

$text = $wsTeste.Cells.Item(17, 1).Value2
$wsTeste.Cells.Item(206, 1).Value = $text
$wsTeste.Cells.Item(206, 2).Value = 1
$wsTeste.Rows.Item(206).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(119, 1).Value2
$wsTeste.Cells.Item(207, 1).Value = $text
$wsTeste.Cells.Item(207, 2).Value = 1
$wsTeste.Rows.Item(207).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(40, 1).Value2
$wsTeste.Cells.Item(208, 1).Value = $text
$wsTeste.Cells.Item(208, 2).Value = 1
$wsTeste.Rows.Item(208).EntireRow.AutoFit()

$text = $wsTreinamento.Cells.Item(298, 1).Value2
$wsTeste.Cells.Item(209, 1).Value = $text
$wsTeste.Cells.Item(209, 2).Value = 1
$wsTeste.Rows.Item(209).EntireRow.AutoFit()

# --- View state: mirror the final selection on each sheet (matches the saved
# workbook state), restoring "Teste" as the active tab at the end.
$wsTreinamento.Range("A302:B337").Select() | Out-Null

$wsTeste.Activate() | Out-Null
$wsTeste.Range("A197").Select() | Out-Null
